# Fix: Fixed try_print_anvisa_register method and added exception handling
# in Google search.
#
# The registration-check results in column D ("Registro") were incorrectly
# reporting "Sucesso" for several rows; after the fix every row reflects the
# real outcome, "Falha". Update all data rows (2-7) in column D accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = "Falha"
    }
}
